# Update LR-pairs sheet (Cxcl13-Ackr4) with new TPM-derived values.
# Sending/Ligand/Receptor/Target cluster labels for rows 2-4 (FAPs -> MuSCs,
# FAPs -> FAPs, FAPs -> ECs edges of the Cxcl13/Ackr4 pair) stay the same;
# only the recomputed metric columns (E:T) change.
# Rows 5-7 (edges sourced from MuSCs) no longer pass the pipeline's
# thresholds with the refreshed TPM values and are removed, shrinking the
# sheet from A1:T7 down to A1:T4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 2: FAPs -> MuSCs (Cxcl13 -> Ackr4) ----
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1148133333333333
$ws.Range("N2").Value = 0.34444
$ws.Range("O2").Value = 0.03343792635928704
$ws.Range("P2").Value = 0.03343792635928704
$ws.Range("Q2").Value = 0.4603839767688889
$ws.Range("R2").Value = 4.14345579092
$ws.Range("S2").Value = 0.03343792635928704
$ws.Range("T2").Value = 0.03343792635928704

# ---- Row 3: FAPs -> FAPs (Cxcl13 -> Ackr4) ----
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.9249645515654102
$ws.Range("P3").Value = 0.9249645515654102
$ws.Range("S3").Value = 0.9249645515654102
$ws.Range("T3").Value = 0.9249645515654102

# ---- Row 4: FAPs -> ECs (Cxcl13 -> Ackr4) ----
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1428303333333333
$ws.Range("N4").Value = 0.428491
$ws.Range("O4").Value = 0.04159752207530271
$ws.Range("P4").Value = 0.04159752207530271
$ws.Range("Q4").Value = 0.5727278788458888
$ws.Range("R4").Value = 5.154550909613
$ws.Range("S4").Value = 0.04159752207530271
$ws.Range("T4").Value = 0.04159752207530271

# ---- Remove the old rows 5-7 (MuSCs-sourced edges no longer present) ----
$ws.Range("A5:T7").EntireRow.Delete()
